$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.913.51"
$ws.Range("E2").Value = "  +4.35%  "
$ws.Range("D3").Value = "3.674.65"
$ws.Range("E3").Value = "  +10.51%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "241.94"
$ws.Range("E5").Value = "  +4.77%  "
$ws.Range("D6").Value = "644.90"
$ws.Range("E6").Value = "  +5.01%  "
$ws.Range("E7").Value = "  +5.01%  "
$ws.Range("D8").Value = "0.402"
$ws.Range("E8").Value = "  +4.53%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  +5.32%  "
$ws.Range("D11").Value = "3.663.42"
$ws.Range("E11").Value = "  +10.15%  "
$ws.Range("D12").Value = "43.85"
$ws.Range("E12").Value = "  +3.22%  "
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("D14").Value = "6.39"
$ws.Range("E14").Value = "  +4.03%  "
$ws.Range("D15").Value = "4.369.19"
$ws.Range("E15").Value = "  +10.69%  "
$ws.Range("D16").Value = "95.796.24"
$ws.Range("E16").Value = "  +4.39%  "
$ws.Range("D17").Value = "0.0000257"
$ws.Range("E17").Value = "  +5.65%  "
$ws.Range("D18").Value = "3.665.65"
$ws.Range("E18").Value = "  +10.25%  "
$ws.Range("D19").Value = "13.42"
$ws.Range("E19").Value = "  +23.98%  "
$ws.Range("D20").Value = "8.06"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "18.76"
$ws.Range("E21").Value = "  +8.53%  "
$ws.Range("D22").Value = "519.56"
$ws.Range("E22").Value = "  +5.66%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "3.44"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").Value = "0.482"
$ws.Range("E24").Value = "  +9.75%  "
$ws.Range("E25").Value = "  +9.37%  "
$ws.Range("D26").Value = "6.80"
$ws.Range("E26").Value = "  +3.51%  "
$ws.Range("D27").Value = "97.45"
$ws.Range("E27").Value = "  +8.09%  "
$ws.Range("D28").Value = "12.67"
$ws.Range("E28").Value = "  +6.37%  "
$ws.Range("D29").Value = "3.16"
$ws.Range("E29").Value = "  +21.10%  "
$ws.Range("D30").Value = "11.70"
$ws.Range("E30").Value = "  +5.17%  "
$ws.Range("E31").Value = "  +2.52%  "
$ws.Range("D33").Value = "33.23"
$ws.Range("E33").Value = "  +17.52%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("B35").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D35").Value = "0.180"
$ws.Range("E35").Value = "  +4.73%  "
$ws.Range("D36").Value = "0.581"
$ws.Range("E36").Value = "  +10.25%  "
$ws.Range("D37").Value = "562.57"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  +9.43%  "
$ws.Range("D39").Value = "7.87"
$ws.Range("E39").Value = "  +5.93%  "
$ws.Range("D40").Value = "0.969"
$ws.Range("E40").Value = "  +11.69%  "
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D43").Value = "1.76"
$ws.Range("E43").Value = "  +5.06%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "5.82"
$ws.Range("E44").Value = "  +7.68%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0432"
$ws.Range("E45").Value = "  +4.42%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "33.82"
$ws.Range("E47").Value = "  +50.28%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "2.22"
$ws.Range("E48").Value = "  +5.33%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "54.81"
$ws.Range("E49").Value = "  +5.71%  "
$ws.Range("D50").Value = "8.33"
$ws.Range("E50").Value = "  +3.94%  "
$ws.Range("E51").Value = "  -2.69%  "
